$d = $word.ActiveDocument

# --- Step 1: paragraph under "Objetivos" heading (para 6) ---
# becomes the short 5-item "Programa" summary list
$d.Paragraphs.Item(6).Range.Text = "1 - Introdução aos Biomateriais`v2 - Interação tecido - implante`v3 - Técnicas de modificação de superfície`v4 - Técnicas de caracterização biológica`v5 - Aspectos práticos no uso de biomateriais"

# --- Step 2: paragraph under "Programa resumido" heading (para 9) ---
# becomes the long objectives paragraph (moved from para 6)
$d.Paragraphs.Item(9).Range.Text = "A ciência dos biomateriais é uma atividade multidisciplinar que envolve a medicina, as ciências naturais e as engenharias, delimitando duas grandes áreas: a biotecnologia e a bioengenharia. A disciplina Biomateriais visa prover aos estudantes fundamentos básicos da ciência de biomateriais, dar uma perspectiva sobre os principais biomateriais aplicados em algumas áreas da medicina e contribuir para a compreensão das interações célula-material. Dessa forma, contribuir para o desenvolvimento da área e certamente alavancar a formação de recursos humanos associados a um melhor uso da infra-estrutura já existente."

# --- Step 3: paragraph under "Programa" heading (para 11) ---
# becomes the long detailed program list (moved up from the Avaliacao/Metodo body)
$d.Paragraphs.Item(11).Range.Text = "1 - Introdução aos Biomateriais`v  1.1- Conceitos básicos de biomateriais; `v  1.2 - Classes de materiais usados na área biomédica;`v  1.3 - Classificação dos biomateriais quanto à resposta biológica`v2 - Interação tecido  implante:`v  2.1 - Histórico da osteointegração; `v  2.2 - Fisiologia do osso;`v  2.3 - Natureza da ligação osso-implante;`v  2.4 - Aspectos superficiais dos implantes.`v3 - Técnicas de modificação da superfície:`v  3.1 - Técnicas para criar uma superfície bioativa: cerâmicas bioativas e biovidros, recobrimentos com fosfatos de  cálcio como transportador de proteínas ósseas morfogenéticas;`v  3.2 - Técnicas para aumentar a rugosidade superficial: usinagem, ataque ácido, jateamento, aspersão térmica. `v4 - Técnicas de caracterização biológica`v  4.1 - Teste em líquido corporal simulado (SBF)`v  4.2 - Cultura de células (in vitro)`v  4.3  Teste com cobaias (in vivo)`v5 - Aspectos práticos no uso de biomateriais`v  5.1- Técnicas de esterilização`v  5.2  Normas técnicas"

# --- Step 4: paragraph under "Avaliacao" heading (para 13), ListBullet ---
# Three bold labels stay put; only the body text that follows each label changes.
# Work from the end of the paragraph backwards so earlier offsets remain valid.

function FindLabelRange($startPos, $endPos, $label) {
    $rng = $d.Range($startPos, $endPos)
    $found = $rng.Find.Execute($label, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("label not found: " + $label)
    }
    return $rng
}

$p13 = $d.Paragraphs.Item(13).Range
$paraStart = $p13.Start
$paraEnd = $p13.End

$metodoLabel = FindLabelRange $paraStart $paraEnd "Método: "
$criterioLabel = FindLabelRange $metodoLabel.End $paraEnd "Critério: "
$normaLabel = FindLabelRange $criterioLabel.End $paraEnd "Norma de recuperação: "

# Replace bodies back-to-front so previously located positions stay valid.
$normaBodyRange = $d.Range($normaLabel.End, $paraEnd - 1)
$normaBodyRange.Text = "Uma prova escrita (Rec) que será composta á NF para obtenção da média final (MF) pelo seguinte critério: `vMF = (Rec+NF)/2"

$criterioBodyRange = $d.Range($criterioLabel.End, $normaLabel.Start)
$criterioBodyRange.Text = "Serão utilizadas duas notas para compor a nota final sendo: NF=(P1+P2)/2`vP1 e P2 serão avaliações escritas (eventualmente a P2 poderá ser substituída por trabalho apresentado por escrito e oral).`v"

$metodoBodyRange = $d.Range($metodoLabel.End, $criterioLabel.Start)
$metodoBodyRange.Text = "As aulas serão expositivas com auxilio do quadro para anotações e empregando-se recursos audiovisuais.`v"

# --- Step 5: paragraph under "Bibliografia" heading (para 15) ---
# becomes "2166002 - Sandra Giacomin Schneider" (moved from the old para 9)
$d.Paragraphs.Item(15).Range.Text = "2166002 - Sandra Giacomin Schneider"
